$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: "Падеж" (due date) - push due date out to 03/12/2020, keep as text
$dueDateCells = @("G2", "G3", "G4", "G5")
foreach ($addr in $dueDateCells) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "03/12/2020"
    $cell.Style = "Normal"
}

# Column J: "Количество" (quantity) - corrected quantity value (decimal), stored as text
$qtyCells = @("J2", "J4", "J5")
foreach ($addr in $qtyCells) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "11414.853"
    $cell.Style = "Normal"
}

# Column U: "email" - updated recipient email addresses for invoice delivery
$emailCells = @("U2", "U3", "U4", "U5")
foreach ($addr in $emailCells) {
    $ws.Range($addr).Value = "energy_ee@bdz.bg;vkunova@bdz.bg"
}
